$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, shifting existing rows 43-45 down to 44-46
$ws.Rows.Item(43).Insert()

# Fill in the new row 43 with the new data
$ws.Cells.Item(43, 1).Value = 8
$ws.Cells.Item(43, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(43, 3).Value = "Coquimbo"
$ws.Cells.Item(43, 4).Value = 45223
$ws.Cells.Item(43, 4).NumberFormat = $ws.Cells.Item(44, 4).NumberFormat
$ws.Cells.Item(43, 5).Value = 4
$ws.Cells.Item(43, 6).Value = 100112013
$ws.Cells.Item(43, 7).Value = "Alcachofa"
$ws.Cells.Item(43, 8).Value = "Española"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 440
$ws.Cells.Item(43, 11).Value = 9000
$ws.Cells.Item(43, 12).Value = 10000
$ws.Cells.Item(43, 13).Value = 9500
$ws.Cells.Item(43, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(43, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(43, 16).Value = 317
$ws.Cells.Item(43, 17).Value = 30
$ws.Cells.Item(43, 18).Value = "Hortaliza"
